# The workbook's "CasosColombia" sheet had an extra (last) data row (row 189)
# appended, which duplicated the figures from the prior day instead of being
# new data; this edit removes that erroneous last row, restoring the sheet to
# rows 1:188 and the view scroll/selection to the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last row (row 189) entirely, shifting nothing below it (it is
# the last row), which also drops the now-unused "####" shared string that
# only that row referenced.
$ws.Rows(189).Delete()

# Reset the view back to the top-left of the data so the saved sheet view
# matches a freshly scrolled-to-top state instead of pointing at the
# now-removed last row.
$ws.Range("B1").Select()
